# Added function to read index in Excel file from data passed in from config.json file
#
# This script reproduces the cell edits made to Sheet1:
#   - C40: "test"
#   - C3:  "h"
#   - D29: "hry"
#   - D31: "hry"
#   - C26: "hey"
# and moves the active selection to F30 (also clears the frozen/scrolled
# "topLeftCell" that was previously set to A10).
#
# New values are written in the same order the unique strings were first
# introduced so that the shared-strings table is rebuilt in the same order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C40").Value = "test"
$ws.Range("C3").Value = "h"
$ws.Range("D29").Value = "hry"
$ws.Range("D31").Value = "hry"
$ws.Range("C26").Value = "hey"

# Update the visible selection to F30 (this also resets any custom
# topLeftCell scroll position on the sheet view).
$ws.Range("F30").Select() | Out-Null
